$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the pokemon name cells (A4:A6 -> pikachu, chikorita, bulbasaur)
# so only totodile/charmander remain in column A, and no sheets are created
# per pokemon (per commit message, the per-pokemon rows are no longer needed here).
$ws.Range("A4:A6").ClearContents()

# Move selection to A4 to match the saved state after the edit.
$ws.Range("A4").Select()
